$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" (column D) values parse as plain numbers (e.g. "216.11").
# Force those specific cells to Text format first so Excel keeps the exact
# string representation (matching multi-thousand values like "26.762.89"
# that already are not auto-numeric and must stay textual too).
$textPriceRows = @(5,10,11,13,16,19,21,22,24,25,28,29,30,32,36,37,38,41,43,45,46,49,51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.762.89"
$ws.Range("E2").Value = "  +0.46%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.650.80"
$ws.Range("E3").Value = "  +0.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.16%  "

# Row 5 - BNB
$ws.Range("D5").Value = "216.11"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.30%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.27%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.08%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.46"
$ws.Range("E10").Value = "  +1.50%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.71%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.879.96"
$ws.Range("E12").Value = "  +0.77%  "

# Row 13 - now Polkadot (was WrappedEther)
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.24"
$ws.Range("E13").Value = "  +3.50%  "

# Row 14 - now WrappedEther (was Polkadot)
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.654.48"
$ws.Range("E14").Value = "  +1.44%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.71%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "66.90"
$ws.Range("E16").Value = "  +5.59%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.808.08"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +1.82%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "221.44"
$ws.Range("E19").Value = "  +1.25%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.13%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "4.43"
$ws.Range("E21").Value = "  +2.79%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "6.36"
$ws.Range("E22").Value = "  +2.50%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +1.16%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.17"
$ws.Range("E24").Value = "  +12.93%  "

# Row 25 - Monero
$ws.Range("D25").Value = "146.99"
$ws.Range("E25").Value = "  -1.27%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.36%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.60%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "7.14"
$ws.Range("E28").Value = "  +4.28%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.98"
$ws.Range("E29").Value = "  +3.58%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0523"
$ws.Range("E30").Value = "  +1.97%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.94%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.45"
$ws.Range("E32").Value = "  +4.10%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +4.37%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +4.09%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.295.25"
$ws.Range("E35").Value = "  +8.23%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.0183"
$ws.Range("E36").Value = "  +5.43%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.41"
$ws.Range("E37").Value = "  +1.29%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "0.836"
$ws.Range("E38").Value = "  +3.49%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +5.06%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  +0.20%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.817"
$ws.Range("E41").Value = "  +3.00%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  -2.90%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "5.47"
$ws.Range("E43").Value = "  +1.17%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.791.83"
$ws.Range("E44").Value = "  +1.10%  "

# Row 45 - Quant
$ws.Range("D45").Value = "93.87"
$ws.Range("E45").Value = "  +1.86%  "

# Row 46 - Aave
$ws.Range("D46").Value = "60.02"
$ws.Range("E46").Value = "  +9.60%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +5.69%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +0.95%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "7.85"
$ws.Range("E49").Value = "  +2.84%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  +3.75%  "

# Row 51 - Mantle
$ws.Range("D51").Value = "0.408"
$ws.Range("E51").Value = "  -0.48%  "
